# Weekly Fruit/Vegetable update: prepend a new week of "Chirimoya" price
# records (Mercado Mayorista Lo Valledor de Santiago, origin "Provincia de
# Limarí", date serial 44463) ahead of the existing history, shifting the
# previously-existing rows 56-95 down to rows 60-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at row 56 (rows 56-95 shift down to rows 60-99).
$ws.Range("A56:A59").EntireRow.Insert()

# New data for the week of 2021-09-24 (serial 44463).
$newRows = @(
    @{Row=56; Calidad="Especial";              Volumen=200; Precio=3000},
    @{Row=57; Calidad="Extra (doble especial)"; Volumen=150; Precio=3300},
    @{Row=58; Calidad="Primera";                Volumen=250; Precio=2700},
    @{Row=59; Calidad="Segunda";                Volumen=200; Precio=2200}
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value  = 6
    $ws.Cells.Item($r, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value  = "Metropolitana"
    $ws.Cells.Item($r, 4).Value  = 44463
    $ws.Cells.Item($r, 5).Value  = 13
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100107
    $ws.Cells.Item($r, 8).Value  = "Otros"
    $ws.Cells.Item($r, 9).Value  = 100107002
    $ws.Cells.Item($r, 10).Value = "Chirimoya"
    $ws.Cells.Item($r, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 12).Value = $nr.Calidad
    $ws.Cells.Item($r, 13).Value = $nr.Volumen
    $ws.Cells.Item($r, 14).Value = $nr.Precio
    $ws.Cells.Item($r, 15).Value = $nr.Precio
    $ws.Cells.Item($r, 16).Value = $nr.Precio
    $ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 15 kilos)"
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 19).Value = $nr.Precio
    $ws.Cells.Item($r, 20).Value = 1
}
